$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 370, pushing existing rows 370-478 down to 372-480.
$ws.Rows.Item(370).Insert()
$ws.Rows.Item(371).Insert()

# Populate the first new row (370) with the new weekly price observation (Primera).
$ws.Range("A370").Value = 5
$ws.Range("B370").Value = "Macroferia Regional de Talca"
$ws.Range("C370").Value = "Maule"
$ws.Range("D370").Value = 45135
$ws.Range("E370").Value = 7
$ws.Range("F370").Value = "Fruta"
$ws.Range("G370").Value = 100101
$ws.Range("H370").Value = "Berries"
$ws.Range("I370").Value = 100101007
$ws.Range("J370").Value = "Kiwi"
$ws.Range("K370").Value = "Hayward"
$ws.Range("L370").Value = "Primera"
$ws.Range("M370").Value = 360
$ws.Range("N370").Value = 12000
$ws.Range("O370").Value = 12000
$ws.Range("P370").Value = 12000
$ws.Range("Q370").Value = "$/bandeja 18 kilos"
$ws.Range("R370").Value = "Provincia de Curicó"
$ws.Range("S370").Value = 667
$ws.Range("T370").Value = 18

# Populate the second new row (371) with the new weekly price observation (Segunda).
$ws.Range("A371").Value = 5
$ws.Range("B371").Value = "Macroferia Regional de Talca"
$ws.Range("C371").Value = "Maule"
$ws.Range("D371").Value = 45135
$ws.Range("E371").Value = 7
$ws.Range("F371").Value = "Fruta"
$ws.Range("G371").Value = 100101
$ws.Range("H371").Value = "Berries"
$ws.Range("I371").Value = 100101007
$ws.Range("J371").Value = "Kiwi"
$ws.Range("K371").Value = "Hayward"
$ws.Range("L371").Value = "Segunda"
$ws.Range("M371").Value = 240
$ws.Range("N371").Value = 10000
$ws.Range("O371").Value = 10000
$ws.Range("P371").Value = 10000
$ws.Range("Q371").Value = "$/bandeja 18 kilos"
$ws.Range("R371").Value = "Provincia de Curicó"
$ws.Range("S371").Value = 556
$ws.Range("T371").Value = 18
